$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "Islands with Magnis"
$ws.Range("B1").Value = "Islands without Magnis"

# Column B updates (rows 4-9)
$ws.Range("B4").Value = "Espanola"
$ws.Range("B5").Value = "San Cristobal"
$ws.Range("B6").Value = "Santa Fe"
$ws.Range("B7").Value = "Champion"
$ws.Range("B8").Value = "Baltra"
$ws.Range("B9").Value = "Enderby"

# Clear B10 (previously "Daphne Major <1983")
$ws.Range("B10").ClearContents()

# Delete row 16 (A16 = "Espanola", now relocated to B4)
$ws.Rows(16).Delete()

# Column widths (stored "width" = ColumnWidth + 5/6, so back the offset out
# to land as close as possible to the target stored widths of 18.5703125 /
# 21.7109375)
$ws.Columns(1).ColumnWidth = 17.736979166666668
$ws.Columns(2).ColumnWidth = 20.877604166666668

# Selection
$ws.Range("A3:A15").Select()
